$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.244.49'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '2.368.47'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'547.78"
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").Value = "'133.76"
$ws.Range("E6").Value = '  -0.82%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +5.18%  '

$ws.Range("E9").Value = '  +3.72%  '

$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("E12").Value = '  -1.11%  '

$ws.Range("D13").Value = "'24.19"
$ws.Range("E13").Value = '  +2.58%  '

$ws.Range("D14").Value = '2.789.34'
$ws.Range("E14").Value = '  +0.91%  '

$ws.Range("D15").Value = '58.171.92'
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").Value = "'0.0000136"
$ws.Range("E16").Value = '  +2.12%  '

$ws.Range("D17").Value = '2.381.98'
$ws.Range("E17").Value = '  +1.38%  '

$ws.Range("D18").Value = "'10.98"
$ws.Range("E18").Value = '  +3.37%  '

$ws.Range("E19").Value = '  +2.73%  '

$ws.Range("D20").Value = "'331.37"
$ws.Range("E20").Value = '  -1.13%  '

$ws.Range("D21").Value = "'6.90"
$ws.Range("E21").Value = '  +2.78%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").Value = "'63.48"
$ws.Range("E23").Value = '  +2.87%  '

$ws.Range("E24").Value = '  -0.92%  '

$ws.Range("E25").Value = '  -0.12%  '

$ws.Range("E26").Value = '  -2.40%  '

$ws.Range("E27").Value = '  -5.94%  '

$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("D29").Value = "'170.47"
$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +1.72%  '

$ws.Range("E31").Value = '  +0.50%  '

$ws.Range("D32").Value = "'18.47"
$ws.Range("E32").Value = '  -0.06%  '

$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("E34").Value = '  -3.92%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").Value = "'4.19"
$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("E37").Value = '  -1.41%  '

$ws.Range("D38").Value = "'1.60"
$ws.Range("E38").Value = '  -2.04%  '

$ws.Range("D39").Value = "'0.413"
$ws.Range("E39").Value = '  +9.08%  '

$ws.Range("D40").Value = "'142.82"
$ws.Range("E40").Value = '  -4.04%  '

$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("D42").Value = "'288.77"

$ws.Range("E43").Value = '  +2.66%  '

$ws.Range("D44").Value = "'0.0518"
$ws.Range("E44").Value = '  +2.57%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = "'0.566"
$ws.Range("E45").Value = '  +0.72%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = "'18.92"
$ws.Range("E46").Value = '  -1.84%  '

$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("D49").Value = "'11.09"
$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("E50").Value = '  +0.77%  '

$ws.Range("E51").Value = '  +0.15%  '
